$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "trainingaudio/16_kotapi2.wav"
$ws.Range("B2").Value = "pngimages/16_icecream.png"

$ws.Range("A3").Value = "trainingaudio/02_pitito3.wav"
$ws.Range("B3").Value = "pngimages/02_pallet.png"

$ws.Range("A4").Value = "trainingaudio/26_kapako1.wav"
$ws.Range("B4").Value = "pngimages/26_pineapple.png"

$ws.Range("A5").Value = "trainingaudio/22_kakoki1.wav"
$ws.Range("B5").Value = "pngimages/22_egg.png"

$ws.Range("A6").Value = "trainingaudio/14_pokoto1.wav"
$ws.Range("B6").Value = "pngimages/14_coffee.png"

$ws.Range("A7").Value = "trainingaudio/05_titopo2.wav"
$ws.Range("B7").Value = "pngimages/05_megaphone.png"

$ws.Range("A8").Value = "trainingaudio/23_patoko1.wav"
$ws.Range("B8").Value = "pngimages/23_lemon.png"

$ws.Range("A9").Value = "trainingaudio/27_pakapa1.wav"
$ws.Range("B9").Value = "pngimages/27_kiwi.png"

$ws.Range("A10").Value = "trainingaudio/12_pokika3.wav"
$ws.Range("B10").Value = "pngimages/12_pie.png"
